# Update "Inscricoes" worksheet with revised enrolment figures
# (Inscritos / Pagos / Inscrições homologadas columns E, F, H)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 5
$ws.Range("E5").Value = 108

# Row 6
$ws.Range("E6").Value = 37

# Row 10
$ws.Range("E10").Value = 369
$ws.Range("F10").Value = 166
$ws.Range("H10").Value = 166

# Row 11
$ws.Range("E11").Value = 246

# Row 12
$ws.Range("E12").Value = 362
$ws.Range("F12").Value = 198
$ws.Range("H12").Value = 198

# Row 15
$ws.Range("E15").Value = 120

# Row 16
$ws.Range("E16").Value = 160

# Row 17
$ws.Range("E17").Value = 69

# Row 21
$ws.Range("E21").Value = 116

# Row 23
$ws.Range("E23").Value = 159
$ws.Range("F23").Value = 69
$ws.Range("H23").Value = 69

# Row 24
$ws.Range("E24").Value = 161
$ws.Range("F24").Value = 78
$ws.Range("H24").Value = 78

# Row 26
$ws.Range("E26").Value = 109
$ws.Range("F26").Value = 62
$ws.Range("H26").Value = 62

# Row 27
$ws.Range("E27").Value = 248
$ws.Range("F27").Value = 118
$ws.Range("H27").Value = 118

# Row 29
$ws.Range("E29").Value = 138
$ws.Range("F29").Value = 77
$ws.Range("H29").Value = 77

# Row 30
$ws.Range("E30").Value = 163
$ws.Range("F30").Value = 90
$ws.Range("H30").Value = 90

# Row 31
$ws.Range("E31").Value = 64

# Row 32
$ws.Range("E32").Value = 148
$ws.Range("F32").Value = 80
$ws.Range("H32").Value = 80

# Row 34
$ws.Range("E34").Value = 166

# Row 35
$ws.Range("E35").Value = 108
$ws.Range("F35").Value = 64
$ws.Range("H35").Value = 64

# Row 40
$ws.Range("E40").Value = 207

# Row 41
$ws.Range("E41").Value = 302
$ws.Range("F41").Value = 128
$ws.Range("H41").Value = 128

# Row 42
$ws.Range("E42").Value = 269
$ws.Range("F42").Value = 141
$ws.Range("H42").Value = 141

# Row 45
$ws.Range("E45").Value = 107
$ws.Range("F45").Value = 46
$ws.Range("H45").Value = 46

# Row 47
$ws.Range("E47").Value = 343

# Row 48
$ws.Range("E48").Value = 155

# Row 50
$ws.Range("E50").Value = 192

# Row 51
$ws.Range("E51").Value = 186
